# Updates cryptos list values (Price / Volume(1h) columns, plus two
# swapped rows) to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.273.00"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3
$ws.Range("D3").Value = "1.866.18"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'242.50"
$ws.Range("E5").Value = "  +3.11%  "

# Row 6
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("D7").Value = "'0.4715"
$ws.Range("E7").Value = "  +0.43%  "

# Row 8
$ws.Range("D8").Value = "'42.76"
$ws.Range("E8").Value = "  -2.61%  "

# Row 9
$ws.Range("D9").Value = "'0.2866"
$ws.Range("E9").Value = "  -0.27%  "

# Row 10
$ws.Range("D10").Value = "'0.06466"
$ws.Range("E10").Value = "  -1.89%  "

# Row 11
$ws.Range("D11").Value = "'20.81"
$ws.Range("E11").Value = "  -4.12%  "

# Row 12
$ws.Range("D12").Value = "'0.07717"
$ws.Range("E12").Value = "  -2.64%  "

# Row 13
$ws.Range("D13").Value = "1.877.87"
$ws.Range("E13").Value = "  +0.29%  "

# Row 14
$ws.Range("D14").Value = "'94.96"
$ws.Range("E14").Value = "  -1.88%  "

# Row 15
$ws.Range("D15").Value = "'0.7062"
$ws.Range("E15").Value = "  +1.93%  "

# Row 16
$ws.Range("D16").Value = "'5.083"
$ws.Range("E16").Value = "  -0.59%  "

# Row 17
$ws.Range("D17").Value = "'269.67"
$ws.Range("E17").Value = "  +0.18%  "

# Row 18
$ws.Range("D18").Value = "30.264.64"
$ws.Range("E18").Value = "  -0.09%  "

# Row 19
$ws.Range("D19").Value = "'13.32"
$ws.Range("E19").Value = "  -5.00%  "

# Row 20
$ws.Range("D20").Value = "'0.000007537"
$ws.Range("E20").Value = "  -2.12%  "

# Row 21
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$ws.Range("D22").Value = "2.110.87"
$ws.Range("E22").Value = "  -0.60%  "

# Row 23
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("D24").Value = "'5.192"
$ws.Range("E24").Value = "  -1.34%  "

# Row 25
$ws.Range("D25").Value = "'6.118"
$ws.Range("E25").Value = "  -1.55%  "

# Row 26
$ws.Range("D26").Value = "'9.328"
$ws.Range("E26").Value = "  -0.89%  "

# Row 27
$ws.Range("D27").Value = "'165.42"
$ws.Range("E27").Value = "  -1.32%  "

# Row 28
$ws.Range("D28").Value = "'18.83"
$ws.Range("E28").Value = "  -0.45%  "

# Row 29
$ws.Range("D29").Value = "'1.911"
$ws.Range("E29").Value = "  -2.04%  "

# Row 30
$ws.Range("D30").Value = "'1.378"
$ws.Range("E30").Value = "  +1.37%  "

# Row 31
$ws.Range("D31").Value = "'0.09846"
$ws.Range("E31").Value = "  -0.49%  "

# Row 32
$ws.Range("E32").Value = "  +2.83%  "

# Row 33
$ws.Range("D33").Value = "'4.244"
$ws.Range("E33").Value = "  -2.75%  "

# Row 34
$ws.Range("D34").Value = "'4.012"
$ws.Range("E34").Value = "  -1.55%  "

# Row 35
$ws.Range("D35").Value = "'0.04738"
$ws.Range("E35").Value = "  -0.48%  "

# Row 36
$ws.Range("D36").Value = "'1.118"
$ws.Range("E36").Value = "  -1.81%  "

# Row 37
$ws.Range("D37").Value = "'0.6901"
$ws.Range("E37").Value = "  -2.05%  "

# Row 38
$ws.Range("D38").Value = "'2.702"
$ws.Range("E38").Value = "  -0.81%  "

# Row 39
$ws.Range("D39").Value = "'0.01842"
$ws.Range("E39").Value = "  -1.81%  "

# Row 40
$ws.Range("D40").Value = "'2.732"
$ws.Range("E40").Value = "  -2.54%  "

# Row 41
$ws.Range("D41").Value = "'6.326"
$ws.Range("E41").Value = "  +1.59%  "

# Row 42
$ws.Range("D42").Value = "'70.34"
$ws.Range("E42").Value = "  -3.75%  "

# Row 43
$ws.Range("D43").Value = "'0.8407"
$ws.Range("E43").Value = "  -0.22%  "

# Row 44
$ws.Range("D44").Value = "'0.9991"
$ws.Range("E44").Value = "  -0.13%  "

# Row 45
$ws.Range("D45").Value = "'1.893"
$ws.Range("E45").Value = "  -3.37%  "

# Row 46
$ws.Range("D46").Value = "'102.22"
$ws.Range("E46").Value = "  -0.51%  "

# Row 47
$ws.Range("D47").Value = "'0.4065"
$ws.Range("E47").Value = "  -2.86%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.215"
$ws.Range("E48").Value = "  +1.09%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.068"
$ws.Range("E49").Value = "  -1.30%  "

# Row 50
$ws.Range("D50").Value = "'928.82"
$ws.Range("E50").Value = "  -1.53%  "

# Row 51
$ws.Range("D51").Value = "'34.69"
$ws.Range("E51").Value = "  +0.23%  "
